# Update cryptocurrency price/volume data per GitHub Actions refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.002.73"
$ws.Range("E2").Value = "  -1.91%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.238.69"
$ws.Range("E3").Value = "  -1.14%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "579.71"
$ws.Range("E5").Value = "  +0.25%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "173.30"
$ws.Range("E6").Value = "  -3.27%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.633"
$ws.Range("E7").Value = "  +1.02%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.236.11"
$ws.Range("E9").Value = "  -1.28%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.123"
$ws.Range("E10").Value = "  -2.40%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.79"
$ws.Range("E11").Value = "  +1.04%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.391"
$ws.Range("E12").Value = "  -2.61%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.799.83"
$ws.Range("E13").Value = "  -1.28%  "
$ws.Range("E14").Value = "  -3.00%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "65.080.11"
$ws.Range("E15").Value = "  -1.83%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "25.82"
$ws.Range("E16").Value = "  -2.10%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.233.42"
$ws.Range("E17").Value = "  -2.03%  "
$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0000160"
$ws.Range("E18").Value = "  -2.18%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "419.64"
$ws.Range("E19").Value = "  -3.49%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.41"
$ws.Range("E20").Value = "  -1.84%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.87"
$ws.Range("E21").Value = "  -2.41%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.22"
$ws.Range("E22").Value = "  -2.36%  "
$ws.Range("E23").Value = "  +0.05%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "70.81"
$ws.Range("E24").Value = "  -1.43%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.66"
$ws.Range("E25").Value = "  -0.24%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.205"
$ws.Range("E26").Value = "  +3.72%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.497"
$ws.Range("E27").Value = "  -1.41%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0000111"
$ws.Range("E28").Value = "  -1.28%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.18"
$ws.Range("E29").Value = "  +3.38%  "
$ws.Range("E30").Value = "  +0.02%  "
$ws.Range("E31").Value = "  -2.59%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.95"
$ws.Range("E32").Value = "  -1.39%  "
$ws.Range("E33").Value = "  +0.08%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.04"
$ws.Range("E34").Value = "  -2.50%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.46"
$ws.Range("E35").Value = "  -1.75%  "
$ws.Range("E36").Value = "  -1.61%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "156.97"
$ws.Range("E37").Value = "  -0.24%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.40"
$ws.Range("E38").Value = "  -1.31%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.840.12"
$ws.Range("E39").Value = "  +2.66%  "
$ws.Range("E40").Value = "  -1.89%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "25.56"
$ws.Range("E41").Value = "  -4.06%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.25"
$ws.Range("E42").Value = "  -0.80%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.727"
$ws.Range("E43").Value = "  -5.90%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "39.59"
$ws.Range("E44").Value = "  -1.70%  "
$ws.Range("E45").Value = "  -4.04%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0634"
$ws.Range("E46").Value = "  -3.56%  "
$ws.Range("B47").Value = "dogwifhat"
$ws.Range("C47").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.19"
$ws.Range("E47").Value = "  -4.12%  "
$ws.Range("B48").Value = "Bittensor"
$ws.Range("C48").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "304.70"
$ws.Range("E48").Value = "  -4.65%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "22.19"
$ws.Range("E49").Value = "  -4.14%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0265"
$ws.Range("E50").Value = "  -0.59%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.102"
$ws.Range("E51").Value = "  -0.46%  "
